# Updates the "ResumoInscricoes" data values as per the "Data update using git" commit.
# Only numeric value cells in columns E (Inscritos), F (Pagos) and H (Inscrições homologadas)
# change; column G (Isenções deferidas) is untouched in this diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "E4"  = 25
    "F4"  = 7
    "H4"  = 21
    "E5"  = 5
    "E6"  = 7
    "F9"  = 12
    "H9"  = 21
    "E17" = 125
    "F17" = 62
    "H17" = 94
    "E18" = 116
    "E19" = 59
    "F19" = 30
    "H19" = 43
    "F28" = 15
    "H28" = 17
    "E33" = 44
    "E34" = 23
    "E36" = 106
    "E37" = 56
    "F37" = 32
    "H37" = 44
    "E38" = 80
    "E40" = 22
    "F46" = 10
    "H46" = 19
    "E47" = 62
    "E48" = 33
    "F50" = 9
    "H50" = 17
    "E61" = 30
    "F61" = 11
    "H61" = 21
    "E62" = 48
    "F62" = 12
    "H62" = 26
    "E67" = 42
    "E71" = 41
    "F71" = 20
    "H71" = 30
    "E77" = 54
    "F77" = 22
    "H77" = 39
    "E81" = 19
    "F81" = 12
    "H81" = 17
    "E89" = 43
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
